# Rename the "Valor médio unitário na extração vegetal" variable label
# (column B) to "Preço médio recebido na extração vegetal" everywhere it
# appears in the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "Valor médio unitário na extração vegetal"
$newText = "Preço médio recebido na extração vegetal"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Text -eq $oldText) {
        $cell.Value = $newText
    }
}
